$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 29.9100170135498
$ws.Range("D2").Value = 0.1000170135498024
$ws.Range("E2").Value = 0.01000340299942136
$ws.Range("C3").Value = 29.88937187194824
$ws.Range("D3").Value = -0.03062812805175952
$ws.Range("E3").Value = 0.0009380822279549783
$ws.Range("C4").Value = 29.96999931335449
$ws.Range("D4").Value = -0.01000068664551179
$ws.Range("E4").Value = 0.0001000137333817179
$ws.Range("B5").Value = 30.03999999999999
$ws.Range("C5").Value = 30.16000175476074
$ws.Range("D5").Value = 0.1200017547607501
$ws.Range("E5").Value = 0.01440042114565922
$ws.Range("B6").Value = 30.21000000000001
$ws.Range("C6").Value = 30.22909355163574
$ws.Range("D6").Value = 0.01909355163573423
$ws.Range("E6").Value = 0.0003645637140664493
$ws.Range("C7").Value = 30.33775901794434
$ws.Range("D7").Value = 0.1177590179443371
$ws.Range("E7").Value = 0.0138671863072147
$ws.Range("C8").Value = 30.40864372253418
$ws.Range("D8").Value = 0.02864372253418423
$ws.Range("E8").Value = 0.0008204628406153337
$ws.Range("C9").Value = 30.61941337585449
$ws.Range("D9").Value = 0.1794133758544945
$ws.Range("E9").Value = 0.0321891594355061
$ws.Range("C10").Value = 30.43139457702637
$ws.Range("D10").Value = -0.04860542297363679
$ws.Range("E10").Value = 0.002362487142446139
$ws.Range("C11").Value = 30.40582847595215
$ws.Range("D11").Value = -0.2841715240478493
$ws.Range("E11").Value = 0.08075345507967739
$ws.Range("C12").Value = 30.52097129821777
$ws.Range("D12").Value = -0.2290287017822266
$ws.Range("E12").Value = 0.05245414624005207
$ws.Range("C13").Value = 30.62080192565918
$ws.Range("D13").Value = -0.319198074340818
$ws.Range("E13").Value = 0.1018874106628864
$ws.Range("C14").Value = 30.70599174499512
$ws.Range("D14").Value = -0.2440082550048857
$ws.Range("E14").Value = 0.05954002851052931
$ws.Range("C15").Value = 31.09451866149902
$ws.Range("D15").Value = 0.07451866149902742
$ws.Range("E15").Value = 0.005553030911606631
$ws.Range("C16").Value = 31.2536449432373
$ws.Range("D16").Value = 0.1336449432373001
$ws.Range("E16").Value = 0.01786097085290118
$ws.Range("C17").Value = 31.33582878112793
$ws.Range("D17").Value = 0.05582878112792855
$ws.Range("E17").Value = 0.003116852802230151
$ws.Range("C18").Value = 31.1854305267334
$ws.Range("D18").Value = -0.194569473266597
$ws.Range("E18").Value = 0.03785727992724101
$ws.Range("C19").Value = 31.38635063171387
$ws.Range("D19").Value = -0.1936493682861311
$ws.Range("E19").Value = 0.03750007783761764
$ws.Range("B20").Value = 31.65000000000001
$ws.Range("C20").Value = 31.94420433044434
$ws.Range("D20").Value = 0.2942043304443303
$ws.Range("E20").Value = 0.08655618805219667
$ws.Range("C21").Value = 32.52571487426758
$ws.Range("D21").Value = 0.6457148742675827
$ws.Range("E21").Value = 0.4169476988504001
$ws.Range("C22").Value = 32.49223709106445
$ws.Range("D22").Value = 0.212237091064452
$ws.Range("E22").Value = 0.04504458282350048
$ws.Range("C23").Value = 32.56542587280273
$ws.Range("D23").Value = 0.1154258728027315
$ws.Range("E23").Value = 0.01332313211227236
$ws.Range("B24").Value = 32.84999999999999
$ws.Range("C24").Value = 32.8115119934082
$ws.Range("D24").Value = -0.03848800659179119
$ws.Range("E24").Value = 0.001481326651409762
$ws.Range("B25").Value = 32.90000000000001
$ws.Range("C25").Value = 32.99087142944336
$ws.Range("D25").Value = 0.09087142944335369
$ws.Range("E25").Value = 0.008257616689078408
$ws.Range("B26").Value = 33.09999999999999
$ws.Range("C26").Value = 32.95927810668945
$ws.Range("D26").Value = -0.1407218933105412
$ws.Range("E26").Value = 0.01980265125690334
$ws.Range("B27").Value = 33.40000000000001
$ws.Range("C27").Value = 33.53360366821289
$ws.Range("D27").Value = 0.1336036682128849
$ws.Range("E27").Value = 0.01784994015993864
$ws.Range("C28").Value = 33.62392044067383
$ws.Range("D28").Value = -0.07607955932617472
$ws.Range("E28").Value = 0.005788099347264939
$ws.Range("B29").Value = 34.09999999999999
$ws.Range("C29").Value = 33.83574295043945
$ws.Range("D29").Value = -0.2642570495605412
$ws.Range("E29").Value = 0.06983178824244232
$ws.Range("B30").Value = 34.40000000000001
$ws.Range("C30").Value = 34.44326400756836
$ws.Range("D30").Value = 0.04326400756835369
$ws.Range("E30").Value = 0.001871774350874565
$ws.Range("B31").Value = 34.90000000000001
$ws.Range("C31").Value = 35.08557510375977
$ws.Range("D31").Value = 0.1855751037597599
$ws.Range("E31").Value = 0.03443811913544566
$ws.Range("C32").Value = 35.61227798461914
$ws.Range("D32").Value = 0.3122779846191435
$ws.Range("E32").Value = 0.097517539677794
$ws.Range("C33").Value = 35.89323043823242
$ws.Range("D33").Value = 0.193230438232419
$ws.Range("E33").Value = 0.03733800225949271
$ws.Range("C34").Value = 36.06398391723633
$ws.Range("D34").Value = -0.236016082763669
$ws.Range("E34").Value = 0.05570359132310707
$ws.Range("C35").Value = 36.57632064819336
$ws.Range("D35").Value = -0.2236793518066378
$ws.Range("E35").Value = 0.05003245242463764
$ws.Range("C36").Value = 37.20608139038086
$ws.Range("D36").Value = -0.09391860961913778
$ws.Range("E36").Value = 0.008820705232792001
$ws.Range("B37").Value = 37.90000000000001
$ws.Range("C37").Value = 37.85781478881836
$ws.Range("D37").Value = -0.04218521118164631
$ws.Range("E37").Value = 0.001779592042440097
$ws.Range("C38").Value = 38.30781173706055
$ws.Range("D38").Value = -0.1921882629394531
$ws.Range("E38").Value = 0.03693632841168437
$ws.Range("B39").Value = 38.90000000000001
$ws.Range("C39").Value = 39.00927352905273
$ws.Range("D39").Value = 0.1092735290527287
$ws.Range("E39").Value = 0.01194070415163754
$ws.Range("B40").Value = 39.40000000000001
$ws.Range("C40").Value = 39.57085037231445
$ws.Range("D40").Value = 0.1708503723144474
$ws.Range("E40").Value = 0.02918984971998531
$ws.Range("B41").Value = 39.90000000000001
$ws.Range("C41").Value = 39.59535598754883
$ws.Range("D41").Value = -0.3046440124511776
$ws.Range("E41").Value = 0.09280797432235323
$ws.Range("B42").Value = 40.09999999999999
$ws.Range("C42").Value = 40.05672073364258
$ws.Range("D42").Value = -0.04327926635741619
$ws.Range("E42").Value = 0.001873094896436177
$ws.Range("B43").Value = 40.59999999999999
$ws.Range("C43").Value = 40.49784469604492
$ws.Range("D43").Value = -0.1021553039550724
$ws.Range("E43").Value = 0.01043570612615324
$ws.Range("B44").Value = 40.90000000000001
$ws.Range("C44").Value = 40.71113967895508
$ws.Range("D44").Value = -0.1888603210449276
$ws.Range("E44").Value = 0.03566822086519311
$ws.Range("B45").Value = 41.20000000000001
$ws.Range("C45").Value = 41.16640853881836
$ws.Range("D45").Value = -0.03359146118165057
$ws.Range("E45").Value = 0.001128386264318337
$ws.Range("C46").Value = 41.36802291870117
$ws.Range("D46").Value = -0.1319770812988281
$ws.Range("E46").Value = 0.01741794998815749
$ws.Range("C47").Value = 42.05647659301758
$ws.Range("D47").Value = 0.256476593017581
$ws.Range("E47").Value = 0.06578024276590586
$ws.Range("C48").Value = 41.7586784362793
$ws.Range("D48").Value = -0.441321563720706
$ws.Range("E48").Value = 0.1947647226048891
$ws.Range("C49").Value = 43.38093185424805
$ws.Range("D49").Value = 0.680931854248044
$ws.Range("E49").Value = 0.4636681901296795
$ws.Range("C50").Value = 44.19462203979492
$ws.Range("D50").Value = 0.4946220397949261
$ws.Range("E50").Value = 0.2446509622508935
$ws.Range("C51").Value = 43.60408782958984
$ws.Range("D51").Value = -0.5959121704101591
$ws.Range("E51").Value = 0.3551113148429465
$ws.Range("C52").Value = 0.06434516906735155
$ws.Range("E52").Value = 3.005327480343232
$ws.Range("E53").Value = 0.06010654960686464
